$d = $word.ActiveDocument

# Move to the end of the document content (last paragraph) and insert a new
# paragraph break followed by the new text "Vishal is topping".
$endRange = $d.Paragraphs(1).Range
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.Move(1, 1) | Out-Null
$endRange.InsertAfter("Vishal is topping")
